$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update roller_rad value (column G, rows 2-5) from "=1in/4" to "=6.375mm".
# The leading apostrophe forces the text "=6.375mm" to be stored literally
# (as a string) instead of being parsed as a formula.
$ws.Range("G2:G5").Value = "'=6.375mm"
